$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:G11 (ME, MAE, MSE, RMSE, SE, N) following the
# "Added filtering options for the Component Analysis" update.
$data = @(
    @(0.03133494016881249,  0.618674247399176,   0.7992451844515852, 0.8940051367031318, 0.9025266721007571, 50),
    @(0.05905192170129663,  0.7500407276805185,  1.336363517888593,  1.156011902139677,  1.166466736141374,  49),
    @(0.003554835396229139, 0.6262718728555511,  0.7969898074296445, 0.8927428562747756, 0.9021829812466056, 48),
    @(-0.03511572363863584, 0.6635828519735995,  0.8873739543468248, 0.9420052836087623, 0.9517525198146293, 46),
    @(-0.02357502177020739, 0.6102703326053934,  0.7645395854549295, 0.8743795431361198, 0.8837201036899427, 46),
    @(0.01163558867306234,  0.5783626329962743,  0.6903657101012227, 0.8308824887438793, 0.8432949623807621, 34),
    @(0.009654851496814876, 0.5447572585674562,  0.6494921422519612, 0.8059107532797668, 0.8183475071966614, 33),
    @(-0.04525758420911288, 0.5348296361048615,  0.5246369517211453, 0.7243182668697133, 0.7466109761084062, 16),
    @(-0.120014643191906,   0.3396452386262637,  0.2564608965434484, 0.5064196841982432, 0.5186064467643261, 10),
    @(-0.248566666372188,   0.301374533835348,   0.1561541005302556, 0.3951633846021865, 0.3434543508584382, 5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j  # B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
